$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the current SIDEBAR_SUBMENU data value (L2) before it gets
# shifted/overwritten below.
$oldSubmenuValue = $ws.Range("L2").Value()

# Insert a new column at M (13th column). This pushes the previous M
# column (KODE_PARAMETER header/values) one to the right, to N.
$ws.Columns.Item(13).Insert()

# Approximate the column widths Excel computed after the insert: the new
# M column takes on roughly the old "best fit" width that L used to have
# (since the relocated submenu text now lives there), while L gets a
# plain, non-best-fit width for its new, shorter-lived text.
$ws.Columns.Item(13).ColumnWidth = 17.43
$ws.Columns.Item(12).ColumnWidth = 14.14

# New header for the inserted column M, formatted like its neighbour L1.
$ws.Range("M1").Value = "SIDEBAR_SUBMENU_SUBMENU"
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)

# The previous SIDEBAR_SUBMENU data value now represents the
# SIDEBAR_SUBMENU_SUBMENU value, so place it in the new M2 cell, keeping
# the formatting that L2 used to have.
$ws.Range("M2").Value = $oldSubmenuValue
$ws.Range("L2").Copy()
$ws.Range("M2").PasteSpecial(-4122)

# Put the new SIDEBAR_SUBMENU value in L2, matching the quote-prefixed
# style used elsewhere in that row (K2).
$ws.Range("L2").Value = "Setup Kelengkapan Kepesertaan"
$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Update the remembered cell selection to match the saved workbook state.
$ws.Range("M12").Select()
